$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# A new entry (SCRIPT/G01P03A/um2205.ssb) is appended as row 12.
#
# Before this edit, row 11 was the last row in the table and used the
# "plain" row style (no bottom border). Row 12 becomes the new last row
# and should use that plain style, while row 11 becomes the last row of
# its own group and picks up the bottom-border "separator" style that
# rows 4 / 7 / 9 already use between entries.
# ---------------------------------------------------------------------

# 1) Give new row 12 the same formatting row 10 already has (plain style,
#    no separating border) by copying formats only (not values).
$ws.Range("A10:E10").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows.Item(12).RowHeight = 43.2

# 2) Fill in the new row's data. Cells are written in the order
#    English(C) -> filename(A) -> translated(D) -> converted(E) so the
#    shared-string table gets the new entries appended in that exact
#    order (39=English, 40=filename, 41=translated, 42=converted).
$ws.Range("C12").Value = " Getting excited over the planet\'s\nparalysis takes too much effort…"
$ws.Range("A12").Value = "SCRIPT/G01P03A/um2205.ssb"
$ws.Range("D12").Value = " Беспокоиться о планетарном\nпараличе слишком утомительно..."
$ws.Range("E12").Value = " Áåòðïëïéóûòÿ ï ðìàîåóàñîïí\nðàñàìéœå òìéšëïí ôóïíéóåìûîï..."
$ws.Range("B12").Value = 64

# 3) Row 11 becomes the last row of its group -> give it the bottom
#    border "separator" formatting (same as row 9).
$ws.Range("A9:E9").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 4) Scroll the view down to the newly added row and select D14, matching
#    where Excel lands after appending this entry.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D14").Select() | Out-Null
